$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checkout payments")
$ws.Rows("12:12").Insert()
$ws.Cells.Item(12,1).Value2 = "IT shipping method"
$ws.Cells.Item(12,40).Value2 = "Consegna standard (2-6 giorni)"
